# Thai translation pass for "Email T-1 [TEMPLATE] Partner email – if RSVP no"
#
# Uses Find/Replace (wdReplaceAll) for strings that are unique across the
# document, and an "anchor + narrowed Range" pattern for strings that repeat
# (" or ", ", ", ". ") so only the intended occurrence is touched.

$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

# --- Language picker line (top of doc) -------------------------------------
Replace-All "English" "ภาษาอังกฤษ"
# NB: the search text intentionally omits the leading space — a leading
# space right after the hyperlink run boundary makes Find latch onto the
# hyperlink's rPr (color/underline) for the whole replacement run. The
# original leading space is left untouched, so it's also dropped here.
Replace-All "/ Portuguese / French / Thai / Vietnamese / Spanish" "/ ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน"

# --- Brief / Target audience table ------------------------------------------
Replace-All "Brief" "บทย่อ"
Replace-All "An email sent to partners in the target country who have RSVPed no. It will be sent via customer.io" "อีเมล์ส่งถึงคู่ค้าที่อยู่ในประเทศเป้าหมายและได้ตอบว่าไม่เข้าร่วม โดยมันจะถูกส่งผ่านทาง customer.io"
Replace-All "Target audience" "กลุ่มเป้าหมาย"
Replace-All "Invited partners who RSVP no" "พาร์ทเนอร์ที่ได้รับเชิญแต่ตอบว่าไม่เข้าร่วม"

# --- Subject line ------------------------------------------------------------
Replace-All ": Thinking of you at " ": พวกเราคิดถึงคุณจากในงาน "

# --- Hero heading --------------------------------------------------------------
Replace-All "We’ll miss you at the " "เราจะคิดถึงคุณจากที่ในงาน "

# --- Greeting ------------------------------------------------------------------
Replace-All "Dear " "เรียนคุณ "

# ", " right after [PARTNER NAME] -> " " (many other ", " exist elsewhere, so
# anchor on the preceding unique text and narrow the range before replacing).
$rng = $d.Content
$rng.Find.Execute("[PARTNER NAME]") | Out-Null
$narrow = $d.Range($rng.End, $d.Content.End)
$narrow.Find.Execute(", ", $true, $false, $false, $false, $false, `
                      $true, 0, $false, " ", 1) | Out-Null

# --- Body paragraphs -------------------------------------------------------
Replace-All "Thank you for taking the time to respond to our invitation to the upcoming " "ขอขอบคุณที่สละเวลาตอบกลับคำเชิญของเราสำหรับงาน "
Replace-All ". We were really looking forward to seeing you there." " ที่กำลังจะมาถึง พวกเราหวังไว้ว่า จะได้พบเจอคุณที่นั่น"
Replace-All "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up. " "แม้ว่าเราจะผิดหวังที่ไม่สามารถพบคุณได้ แต่เราก็เข้าใจดีว่าปัญหาเกี่ยวกับตารางเวลาที่ขัดแย้งและภาระผูกพันอื่นๆ บางครั้งก็เกิดขึ้นได้ "
Replace-All "If you’re comfortable sharing it with us, we’d like to know why you responded no. Please reply to this email as your feedback could help us make improvements in our event planning processes and better serve you in the future." "หากคุณไม่ขัดข้องที่จะแบ่งปันกับเรา พวกเราก็ต้องการทราบว่า คุณตอบปฏิเสธคำเชิญเพราะอะไร โปรดตอบกลับอีเมล์นี้ เนื่องจากข้อคิดเห็นหรือคำติชมของคุณจะช่วยให้เราได้ปรับปรุงพัฒนากระบวนการวางแผนกิจกรรมของเราและให้บริการคุณได้ดียิ่งขึ้นในอนาคต"
Replace-All "We hope to see you at our future events. " "เราหวังว่า จะได้พบคุณในกิจกรรมของเราในอนาคต "

# --- "contact us via live chat or WhatsApp." paragraph ---------------------
Replace-All "If you have any questions, please contact us via " "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง "
Replace-All "live chat" "แชทสด"

# " or " before the WhatsApp hyperlink -> " หรือทาง " (there are two " or "
# runs in the doc; anchor on the text right before this one).
$rng = $d.Content
$rng.Find.Execute("แชทสด") | Out-Null
$narrow = $d.Range($rng.End, $d.Content.End)
$narrow.Find.Execute(" or ", $true, $false, $false, $false, $false, `
                      $true, 0, $false, " หรือทาง ", 1) | Out-Null

# ". " right after the WhatsApp hyperlink -> " " (anchor on "WhatsApp").
$rng = $d.Content
$rng.Find.Execute("WhatsApp") | Out-Null
$narrow = $d.Range($rng.End, $d.Content.End)
$narrow.Find.Execute(". ", $true, $false, $false, $false, $false, `
                      $true, 0, $false, " ", 1) | Out-Null

# --- "contact your country manager" paragraph -------------------------------
Replace-All "If you have any questions, please contact your country manager, " "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ "
Replace-All ", at " " ที่ "

# second " or " (before [WHATSAPP NO]) -> " หรือ "
$rng = $d.Content
$rng.Find.Execute("[EMAIL ADDRESS]") | Out-Null
$narrow = $d.Range($rng.End, $d.Content.End)
$narrow.Find.Execute(" or ", $true, $false, $false, $false, $false, `
                      $true, 0, $false, " หรือ ", 1) | Out-Null

Replace-All " (WhatsApp). " " (WhatsApp) "

# --- Comment text ------------------------------------------------------------
# The simulated Word OM only writes a comment's body through a direct
# Range.Text assignment (Find/ParagraphFormat writes on Comment.Range are not
# wired through), so use that for the lone comment in the doc.
$d.Comments.Item(1).Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
